$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos-list refresh diff (commit "Updated
# cryptos list ... with GitHub Actions"). The Price/Volume columns (D/E)
# hold numeric-looking text (e.g. "1.000", "29.116.55", "  -1.99%  ") that
# must stay literal text rather than being auto-parsed into a number when
# assigned through .Value. $q is a single apostrophe, prepended to force
# text interpretation (same as typing `'1.000` into Excel); the Style is
# then reset to "Normal" so the quote-prefix flag that trick adds to the
# cell does not linger as a spurious formatting change.
$q = "'"

$ws.Range("D2").Value = $q + '29.116.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = $q + '  -1.99%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = $q + '1.851.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = $q + '  +0.05%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = $q + '0.6955'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = $q + '  -4.46%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = $q + '238.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = $q + '  -0.93%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = $q + '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = $q + '0.07633'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = $q + '  +7.74%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = $q + '0.3029'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = $q + '  -2.94%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = $q + '23.37'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = $q + '  -4.10%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = $q + '0.08126'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = $q + '  -1.24%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = $q + '0.7267'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = $q + '  -2.33%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = $q + '5.219'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = $q + '  -1.62%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = $q + '1.812.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = $q + '  -2.69%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = $q + '89.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = $q + '  -3.41%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = $q + '29.125.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = $q + '  -1.95%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = $q + '5.773'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = $q + '  -3.68%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = $q + '13.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = $q + '  -1.05%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = $q + '0.000007737'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = $q + '  -0.54%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = $q + '236.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = $q + '  -4.66%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = $q + '  -0.02%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = $q + '2.097.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = $q + '  -0.46%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = $q + '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = $q + '  +0.12%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = $q + '7.617'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = $q + '  -1.06%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = $q + '8.980'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = $q + '  -1.89%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = $q + '161.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = $q + '  -1.01%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = $q + '0.1444'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = $q + '  -5.54%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = $q + '18.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = $q + '  -2.44%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = $q + '1.984'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = $q + '  -1.13%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = $q + '1.410'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = $q + '  -1.79%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = $q + '4.480'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = $q + '  -0.58%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = $q + '1.487'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = $q + '  -2.26%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = $q + '4.011'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = $q + '  -4.06%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = $q + '0.05225'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = $q + '  -0.78%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = $q + '1.188'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = $q + '  -3.23%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = $q + '0.7007'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = $q + '  -6.79%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = $q + '1.005'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = $q + '  +0.72%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = $q + '2.652'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = $q + '  -1.54%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = $q + '0.01854'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = $q + '2.680'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = $q + '  -2.04%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = $q + '0.9320'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = $q + '  +7.65%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = $q + '6.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = $q + '  +0.34%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = $q + '1.080.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = $q + '  +3.14%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = $q + '0.4262'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = $q + '  -4.25%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = $q + '70.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = $q + '  -1.09%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = $q + '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = $q + '  +0.00%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = $q + '103.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = $q + '  -0.63%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = $q + '1.776'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = $q + '  -2.11%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = $q + '1.993.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = $q + '  -0.67%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = $q + '9.190'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = $q + '  -3.65%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = $q + '6.997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = $q + '  -6.12%  '
$ws.Range("E51").Style = "Normal"
